$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Work bottom-up so row indices of not-yet-processed rows stay stable.

# --- Row 46 (was "5 <tab> ... <tab> 100.0") -> deleted entirely ---
$t.Rows.Item(46).Delete()

# --- Row 45 (was "95 <tab> ... <tab> 100.0") -> collapses to just "95" ---
$t.Rows.Item(45).Cells.Item(1).Range.Text = "95"

# --- New row inserted after (original) row 44, containing "0.05" ---
$newRow1 = $t.Rows.Add($t.Rows.Item(45))
$newRow1.Cells.Item(1).Range.Text = "0.05"

# --- Row 44 (was "100 <tab> ... <tab> 100.0") -> collapses to just "99.95" ---
$t.Rows.Item(44).Cells.Item(1).Range.Text = "99.95"

# --- New row inserted after (original) row 12, containing "0.05014" ---
$newRow2 = $t.Rows.Add($t.Rows.Item(13))
$newRow2.Cells.Item(1).Range.Text = "0.05014"

# --- Row 12: 0.02017 -> 0.00026 ---
$t.Rows.Item(12).Cells.Item(1).Range.Text = "0.00026"

# --- Row 11: 0.00025 -> 0.00023 ---
$t.Rows.Item(11).Cells.Item(1).Range.Text = "0.00023"

# --- Row 9 (0.00014) -> deleted entirely ---
$t.Rows.Item(9).Delete()

# --- Row 7: 0.00013 -> 0.00014 ---
$t.Rows.Item(7).Cells.Item(1).Range.Text = "0.00014"

# --- Row 6: 0.00048 -> 0.00052 ---
$t.Rows.Item(6).Cells.Item(1).Range.Text = "0.00052"

# --- Row 5: 0.00005 -> 0.00002 ---
$t.Rows.Item(5).Cells.Item(1).Range.Text = "0.00002"

# --- Row 4: 105 -> 305 ---
$t.Rows.Item(4).Cells.Item(1).Range.Text = "305"

# --- Row 3: 95 -> 0M ---
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# --- Row 2: 0.05 -> 0M ---
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"

# --- Row 1: 99.95 -> 0M ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"

Write-Host "Done. Final row count:" $t.Rows.Count
